# Bring the sheet up to date: add a "Calculus" entry next to the existing
# "Evan Miller" row (shifting the 1234 value over to column C), and add a
# second row repeating "Evan Miller"/1234.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: A1 already holds "Evan Miller" / B1 already holds 1234.
# Insert the new "Calculus" label in B1 and push the numeric value to C1.
$ws.Range("B1").Value = "Calculus"
$ws.Range("C1").Value = 1234

# Row 2: new row duplicating the original A1/B1 pair.
$ws.Range("A2").Value = "Evan Miller"
$ws.Range("B2").Value = 1234

# Leave the selection on B3, as in the saved workbook.
$ws.Range("B3").Select()

# Nudge the sheet-tab area ratio back down slightly (cosmetic, matches the
# small adjustment made from the author's "old laptop").
$excel.ActiveWindow.TabRatio = 0.285
